$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45188 to 45189 for all existing data rows (2-491)
for ($r = 2; $r -le 491; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45189
}

# Row 491 gains an explicit row height (customHeight) in the new file
$ws.Rows.Item(491).RowHeight = 15

# Add new row 492
$ws.Range("A492").Value2 = "A 43938-2023"
$ws.Range("B492").Value2 = 45187
$ws.Range("B492").NumberFormat = "YYYY-MM-DD"
$ws.Range("C492").Value2 = 45189
$ws.Range("C492").NumberFormat = "YYYY-MM-DD"
$ws.Range("D492").Value2 = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E492").Value2 = "ÅTVIDABERG"
$ws.Range("F492").Value2 = "Övriga Aktiebolag"
$ws.Range("G492").Value2 = 2
$ws.Range("H492").Value2 = 0
$ws.Range("I492").Value2 = 0
$ws.Range("J492").Value2 = 0
$ws.Range("K492").Value2 = 0
$ws.Range("L492").Value2 = 0
$ws.Range("M492").Value2 = 0
$ws.Range("N492").Value2 = 0
$ws.Range("O492").Value2 = 0
$ws.Range("P492").Value2 = 0
$ws.Range("Q492").Value2 = 0
$ws.Range("R492").WrapText = $true
$ws.Range("R492").Value2 = ""
$ws.Rows.Item(492).RowHeight = 15

# Add new row 493
$ws.Range("A493").Value2 = "A 43940-2023"
$ws.Range("B493").Value2 = 45187
$ws.Range("B493").NumberFormat = "YYYY-MM-DD"
$ws.Range("C493").Value2 = 45189
$ws.Range("C493").NumberFormat = "YYYY-MM-DD"
$ws.Range("D493").Value2 = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E493").Value2 = "ÅTVIDABERG"
$ws.Range("F493").Value2 = "Övriga Aktiebolag"
$ws.Range("G493").Value2 = 3.7
$ws.Range("H493").Value2 = 0
$ws.Range("I493").Value2 = 0
$ws.Range("J493").Value2 = 0
$ws.Range("K493").Value2 = 0
$ws.Range("L493").Value2 = 0
$ws.Range("M493").Value2 = 0
$ws.Range("N493").Value2 = 0
$ws.Range("O493").Value2 = 0
$ws.Range("P493").Value2 = 0
$ws.Range("Q493").Value2 = 0
$ws.Range("R493").WrapText = $true
$ws.Range("R493").Value2 = ""
